$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-like Price cells (column D) to avoid numeric auto-coercion,
# then set value, then reset style so no stray quotePrefix/style lingers.
$priceCells = @(
    "D2",
    "D3",
    "D5",
    "D6",
    "D7",
    "D8",
    "D10",
    "D11",
    "D13",
    "D16",
    "D17",
    "D18",
    "D19",
    "D22",
    "D23",
    "D25",
    "D29",
    "D30",
    "D31",
    "D32",
    "D33",
    "D34",
    "D35",
    "D36",
    "D37",
    "D38",
    "D40",
    "D41",
    "D42",
    "D43",
    "D45",
    "D46",
    "D48",
    "D51"
)
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '69.295.28'
$ws.Range('E2').Value = '  -3.46%  '
$ws.Range('D3').Value = '3.513.35'
$ws.Range('E3').Value = '  -4.85%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '582.68'
$ws.Range('E5').Value = '  -1.25%  '
$ws.Range('D6').Value = '173.88'
$ws.Range('D7').Value = '0.622'
$ws.Range('E7').Value = '  +0.60%  '
$ws.Range('D8').Value = '3.506.31'
$ws.Range('E8').Value = '  -4.83%  '
$ws.Range('D10').Value = '0.190'
$ws.Range('E10').Value = '  -6.53%  '
$ws.Range('D11').Value = '6.73'
$ws.Range('E11').Value = '  +5.29%  '
$ws.Range('E12').Value = '  -3.26%  '
$ws.Range('D13').Value = '47.03'
$ws.Range('E13').Value = '  -6.11%  '
$ws.Range('E15').Value = '  -1.57%  '
$ws.Range('D16').Value = '4.071.60'
$ws.Range('E16').Value = '  -4.94%  '
$ws.Range('D17').Value = '8.71'
$ws.Range('E17').Value = '  -4.04%  '
$ws.Range('D18').Value = '69.249.64'
$ws.Range('E18').Value = '  -3.55%  '
$ws.Range('D19').Value = '3.511.04'
$ws.Range('E19').Value = '  -4.83%  '
$ws.Range('E20').Value = '  -1.42%  '
$ws.Range('E21').Value = '  -4.22%  '
$ws.Range('D22').Value = '11.20'
$ws.Range('E22').Value = '  -4.42%  '
$ws.Range('D23').Value = '0.904'
$ws.Range('E23').Value = '  -4.55%  '
$ws.Range('E24').Value = '  -10.02%  '
$ws.Range('D25').Value = '97.87'
$ws.Range('E25').Value = '  -5.88%  '
$ws.Range('E27').Value = '  -0.59%  '
$ws.Range('E28').Value = '  +0.13%  '
$ws.Range('D29').Value = '2.66'
$ws.Range('E29').Value = '  -7.00%  '
$ws.Range('D30').Value = '9.45'
$ws.Range('E30').Value = '  -7.85%  '
$ws.Range('D31').Value = '32.97'
$ws.Range('E31').Value = '  -6.98%  '
$ws.Range('D32').Value = '8.72'
$ws.Range('E32').Value = '  -6.11%  '
$ws.Range('D33').Value = '3.20'
$ws.Range('E33').Value = '  -8.07%  '
$ws.Range('B34').Value = 'NEARProtocol'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D34').Value = '7.30'
$ws.Range('E34').Value = '  -1.31%  '
$ws.Range('B35').Value = 'Mantle'
$ws.Range('C35').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D35').Value = '1.36'
$ws.Range('E35').Value = '  -6.51%  '
$ws.Range('D36').Value = '596.03'
$ws.Range('E36').Value = '  +5.16%  '
$ws.Range('D37').Value = '3.62'
$ws.Range('E37').Value = '  -15.97%  '
$ws.Range('D38').Value = '10.90'
$ws.Range('E38').Value = '  -3.84%  '
$ws.Range('E39').Value = '  -5.11%  '
$ws.Range('D40').Value = '57.34'
$ws.Range('E40').Value = '  -3.75%  '
$ws.Range('D41').Value = '0.999'
$ws.Range('D42').Value = '0.0439'
$ws.Range('E42').Value = '  -6.20%  '
$ws.Range('D43').Value = '0.336'
$ws.Range('E43').Value = '  -5.16%  '
$ws.Range('E44').Value = '  -6.98%  '
$ws.Range('D45').Value = '3.413.18'
$ws.Range('E45').Value = '  -9.25%  '
$ws.Range('D46').Value = '33.37'
$ws.Range('E46').Value = '  -6.61%  '
$ws.Range('E47').Value = '  -9.16%  '
$ws.Range('D48').Value = '2.91'
$ws.Range('E48').Value = '  +0.06%  '
$ws.Range('E49').Value = '  -7.74%  '
$ws.Range('E50').Value = '  -0.82%  '
$ws.Range('D51').Value = '5.74'
$ws.Range('E51').Value = '  +17.10%  '

foreach ($addr in $priceCells) {
    $ws.Range($addr).Style = "Normal"
}
